$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: B9 keeps its text ("Crear la segunda parte con la consulta de los
# datos de las pelis") but its format is refreshed (re-applying wrap text
# collapses it back onto the shared wrap/vcenter style used elsewhere in
# column B).
$ws.Range("B9").WrapText = $true

# Row 10: mark the "Crear los templates..." task as completed.
$ws.Range("D10").Value = 1

# Row 11: bump priority from Media to Alta.
$ws.Range("C11").Value = "Alta"

# Row 12: new task - "Terminar el readme..." / Alta / note about gifs.
$ws.Range("B12").Value = "Terminar el readme del proyecto y el general"
$ws.Range("C12").Value = "Alta"
$ws.Range("E12").Value = "Hacer los gifs y poner el link"

# Row 13: new task - "Hacer la pagina de error." / Media / note about DIRS
# setting. The three filled-in cells get their font size bumped from 12pt
# to 14pt (same font families as before), and the row grows to fit.
$ws.Range("B13").Value = "Hacer la página de error."
$ws.Range("C13").Value = "Media"
$ws.Range("E13").Value = "HAY QUE CONFIGURAR EL 'DIRS' EN SETTINGS PARA PONER LOS TEMPLATES DE ERROR ONLINE"

$ws.Range("B13").Font.Size = 14
$ws.Range("C13").Font.Size = 14
$ws.Range("E13").Font.Size = 14

$ws.Rows.Item(13).RowHeight = 54

# Selection / scroll position left where the edit happened.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("B13:C13").Select
